$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.312.52"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "1.707.43"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5308"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2655"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06587"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.93%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07640"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.568"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("D13").Value = "1.725.14"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "1.939.78"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5725"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.56%  "
$ws.Range("D16").Value = "0.0₅8167"
$ws.Range("E16").Value = "  -2.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.32%  "
$ws.Range("D18").Value = "27.301.01"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.669"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.960"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.766"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.19%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1215"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.252"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05382"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.291"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.496"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.425"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.424"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9482"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5857"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01628"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.869"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.005"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "1.044.67"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8401"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").Value = "1.847.95"
$ws.Range("E45").Value = "  -1.96%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.004"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.071"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05242"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
